# Auto-generated edit script
# Applies the "Updated symbol list" data refresh described in the diff:
# updates Price (column D) and Volume(1h) (column E) values for many rows,
# and re-shuffles Coin name / Link (columns B / C) for rows 6-17 to reflect
# the new ranking order.
#
# Numeric-looking values (columns D and E) must remain stored as TEXT
# (matching the original t="inlineStr" cell type), so we force the
# NumberFormat to "@" (Text) before assigning them - otherwise Excel's
# COM automation would silently coerce strings like "307.75" or "0.08%"
# into numeric / percentage values and lose exact formatting (trailing
# zeros, "%" suffix, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.08%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.17%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.126"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07617"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.62%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.624"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.37%"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.476"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.15%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9036"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.81%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1095"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.99%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1771"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.08%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09249"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.64%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04258"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.95%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.51%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001248"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.86%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005839"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.23%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.362"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.14%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.250"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.39%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.541"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.45%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.57%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04157"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.68%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001221"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.65%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004085"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.46%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.41%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02415"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "2.60%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05201"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.67%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007762"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.28%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1299"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.95%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006952"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.84%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.11%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008069"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.95%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3055"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006730"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.28%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.21%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03204"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "857.65%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004199"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-16.05%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.21%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.21%"
